# Update the dSF (column F) values to reflect repulled data / recalculated mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -3
    8  = 1
    10 = -4
    17 = -3
    19 = -1
    26 = 6
    31 = 6
    34 = 4
    35 = -2
    37 = 5
    39 = 0
    42 = 3
    44 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
